$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing row (Kim / row 107) - table now ends at row 106
$ws.Rows.Item(107).Delete()

# --- N63:N67: values rotate up by one (old N63 wraps around to N67) ---
$ws.Range("N63").Value = "SALE INVENTORY"
$ws.Range("N64").Value = "DC5-FINANCIAL"
$ws.Range("N65").Value = "MAIN ST STATION, JEFFERSON"
$ws.Range("N66").Value = "1411 S. MAIN STREET"
$ws.Range("N67").Value = "https://maps.app.goo.gl/cmXWoDMNmqSXtfnh7"

# --- I76:I106 / J76:J106: meeting notes block reorganized, numbered list renumbered/shifted ---
$ws.Range("I76").ClearContents()
$ws.Range("I77").ClearContents()
$ws.Range("I78").ClearContents()
$ws.Range("I79").ClearContents()
$ws.Range("I80").ClearContents()
$ws.Range("I81").ClearContents()
$ws.Range("I82").Value = "1)"
$ws.Range("I83").Value = "2)"
$ws.Range("I84").Value = "3)"
$ws.Range("I85").Value = "4)"
$ws.Range("I86").Value = "5)"
$ws.Range("I87").Value = "6)"
$ws.Range("I88").Value = "7)"
$ws.Range("I89").Value = "8)"
$ws.Range("I90").Value = "9)"
$ws.Range("I91").Value = "10)"
$ws.Range("I92").Value = "11)"
$ws.Range("I93").Value = "12)"
$ws.Range("I94").Value = "13)"
$ws.Range("I95").Value = "14)"
$ws.Range("I96").Value = "15)"
$ws.Range("I97").Value = "16)"
$ws.Range("I98").Value = "17)"
$ws.Range("I99").Value = "18)"
$ws.Range("I100").Value = "19)"
$ws.Range("I101").Value = "20)"
$ws.Range("I102").Value = "21)"
$ws.Range("I103").ClearContents()
$ws.Range("I104").ClearContents()
$ws.Range("I105").ClearContents()
$ws.Range("I106").ClearContents()

$ws.Range("J76").Value = "PAID TRAINING/MEETING, `nBADGER EFFICIENCY VIRTUAL WORKSHOP"
$ws.Range("J77").Value = "Please plan to log in 10 mins prior to the start of this meeting"
$ws.Range("J78").Value = "Please have your camera on for this meeting"
$ws.Range("J79").Value = "Please click on the link below to join this meeting"
$ws.Range("J80").Value = "https://meet.google.com/hnc-ibhc-qoh"
$ws.Range("J81").ClearContents()
$ws.Range("J82").Value = "Casey"
$ws.Range("J83").Value = "Elijah"
$ws.Range("J84").Value = "Greg"
$ws.Range("J85").Value = "Heidi"
$ws.Range("J86").Value = "Jake S"
$ws.Range("J87").Value = "Jerry D"
$ws.Range("J88").Value = "Joseph"
$ws.Range("J89").Value = "Josh S"
$ws.Range("J90").Value = "Joshua M"
$ws.Range("J91").Value = "Justin"
$ws.Range("J92").Value = "Kirsten"
$ws.Range("J93").Value = "Lashaun"
$ws.Range("J94").Value = "Makeda"
$ws.Range("J95").Value = "Nick"
$ws.Range("J96").Value = "Robyn"
$ws.Range("J97").Value = "Serena"
$ws.Range("J98").Value = "Sonia"
$ws.Range("J99").Value = "Taya"
$ws.Range("J100").Value = "Taylor"
$ws.Range("J101").Value = "Via"
$ws.Range("J102").Value = "Wyatt"
$ws.Range("J103").ClearContents()
$ws.Range("J104").ClearContents()
$ws.Range("J105").Value = "Office"
$ws.Range("J106").Value = "Kim"
